$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.162.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.277.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.05%  "
$ws.Range("E7").Value = "  +4.05%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -2.42%  "
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.402"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.848.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("E13").Value = "  -3.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "66.165.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.37"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.13%  "
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.266.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "434.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("E19").Value = "  -2.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.61%  "
$ws.Range("E21").Value = "  -3.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.99%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.422.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.68%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.505"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("E26").Value = "  +1.29%  "
$ws.Range("E27").Value = "  -5.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.53%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.61%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E33").Value = "  -2.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.50%  "
$ws.Range("E35").Value = "  -2.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.774.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.55%  "
$ws.Range("E44").Value = "  -2.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0658"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "320.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("E47").Value = "  -2.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.63%  "
$ws.Range("E49").Value = "  -2.25%  "
$ws.Range("E50").Value = "  +2.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.06%  "
